$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20 ("Adjustment Down") switches from an Income/In classification to an
# Expense/Out classification. Copy formatting from an existing Expense/Out
# row (row 2) so the fill/theme-color matches exactly, then set the values.
$ws.Range("B2").Copy()
$ws.Range("B20").PasteSpecial(-4122)
$ws.Range("B20").Value = "Expense"

$ws.Range("C2").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("C20").Value = "Out"

# D19 and D20 both get the new "Manual Adjustment (not sure if needed)" text,
# highlighted with a yellow fill to flag it as uncertain.
$ws.Range("D19").Value = "Manual Adjustment (not sure if needed)"
$ws.Range("D19").Interior.Color = 65535

$ws.Range("D20").Value = "Manual Adjustment (not sure if needed)"
$ws.Range("D20").Interior.Color = 65535

# Leave the selection where the editor ended up working.
$ws.Range("D23").Select() | Out-Null
